$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H104").Value = 119.2
$ws.Range("I104").Value = 119.2
$ws.Range("K104").Value = 357.6
$ws.Range("M104").Value = 1389.4
$ws.Range("H135").Value = 8478159
$ws.Range("I135").Value = 9804929
$ws.Range("K135").Value = 88244361
$ws.Range("M135").Value = -88241826

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3423.7917
$ws.Range("I45").Value = 2831.75
$ws.Range("J45").Value = 4015.8333
$ws.Range("K45").Value = 2831.75
$ws.Range("L45").Value = 4015.8333
$ws.Range("M45").Value = -2454.75
$ws.Range("N45").Value = -4769.8333
$ws.Range("H61").Value = 2539.077
$ws.Range("I61").Value = 1654.3334
$ws.Range("K61").Value = 1654.3334
$ws.Range("M61").Value = -1442.3334
$ws.Range("H97").Value = 5293734.5
$ws.Range("I97").Value = 2000.4445
$ws.Range("J97").Value = 37044136
$ws.Range("K97").Value = 2000.4445
$ws.Range("L97").Value = 37044136
$ws.Range("M97").Value = -1504.4445
$ws.Range("N97").Value = -37045128
$ws.Range("H132").Value = 1390.75
$ws.Range("I132").Value = 1341.4615
$ws.Range("J132").Value = 1604.3334
$ws.Range("K132").Value = 4024.3845
$ws.Range("L132").Value = 4813.0002
$ws.Range("M132").Value = -1494.3845
$ws.Range("N132").Value = -9873.0002
$ws.Range("H136").Value = 2539.077
$ws.Range("I136").Value = 1654.3334
$ws.Range("K136").Value = 4963.0002
$ws.Range("M136").Value = -2413.0002

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1961.5883
$ws.Range("J58").Value = 4382.2
$ws.Range("L58").Value = 4382.2
$ws.Range("N58").Value = -4788.2
$ws.Range("H86").Value = 52341.445
$ws.Range("J86").Value = 39356.332
$ws.Range("L86").Value = 39356.332
$ws.Range("N86").Value = -41602.332
$ws.Range("H89").Value = 52341.445
$ws.Range("J89").Value = 39356.332
$ws.Range("L89").Value = 196781.66
$ws.Range("N89").Value = -208013.66
$ws.Range("H94").Value = 7115.5884
$ws.Range("I94").Value = 14753
$ws.Range("K94").Value = 14753
$ws.Range("M94").Value = -14302
$ws.Range("H134").Value = 3534.139
$ws.Range("I134").Value = 3502.182
$ws.Range("J134").Value = 3584.3572
$ws.Range("K134").Value = 10506.546
$ws.Range("L134").Value = 10753.0716
$ws.Range("M134").Value = -7971.545999999998
$ws.Range("N134").Value = -15823.0716
$ws.Range("H136").Value = 1961.5883
$ws.Range("J136").Value = 4382.2
$ws.Range("L136").Value = 13146.6
$ws.Range("N136").Value = -18246.6

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 25150520
$ws.Range("J4").Value = 746.125
$ws.Range("L4").Value = 2238.375
$ws.Range("N4").Value = -2462.375
$ws.Range("H126").Value = 3449.75
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("H129").Value = 1708.9412
$ws.Range("I129").Value = 557.5
$ws.Range("J129").Value = 2063.2307
$ws.Range("K129").Value = 1672.5
$ws.Range("L129").Value = 6189.6921
$ws.Range("M129").Value = 3327.5
$ws.Range("N129").Value = -16189.6921
$ws.Range("H131").Value = 771691.5600000001
$ws.Range("I131").Value = 1177804.9
$ws.Range("K131").Value = 3533414.7
$ws.Range("M131").Value = -3528374.7
$ws.Range("H134").Value = 3365.75
$ws.Range("I134").Value = 1280.2667
$ws.Range("K134").Value = 3840.800099999999
$ws.Range("M134").Value = 1229.199900000001
$ws.Range("H136").Value = 3845.3845
$ws.Range("I136").Value = 1330
$ws.Range("J136").Value = 4600
$ws.Range("K136").Value = 3990
$ws.Range("L136").Value = 13800
$ws.Range("M136").Value = 1110
$ws.Range("N136").Value = -24000
$ws.Range("H138").Value = 6326.778
$ws.Range("I138").Value = 4988.9
$ws.Range("K138").Value = 14966.7
$ws.Range("M138").Value = -9826.699999999999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 22298950
$ws.Range("I80").Value = 161070.72
$ws.Range("J80").Value = 41669590
$ws.Range("K80").Value = 161070.72
$ws.Range("L80").Value = 41669590
$ws.Range("M80").Value = -160072.72
$ws.Range("N80").Value = -41671586
$ws.Range("H83").Value = 22298950
$ws.Range("I83").Value = 161070.72
$ws.Range("J83").Value = 41669590
$ws.Range("K83").Value = 805353.6
$ws.Range("L83").Value = 208347950
$ws.Range("M83").Value = -800361.6
$ws.Range("N83").Value = -208357934
$ws.Range("H97").Value = 33333712
$ws.Range("I97").Value = 453.8
$ws.Range("J97").Value = 200000000
$ws.Range("K97").Value = 453.8
$ws.Range("L97").Value = 200000000
$ws.Range("M97").Value = 42.19999999999999
$ws.Range("N97").Value = -200000992
$ws.Range("H122").Value = 4561.7095
$ws.Range("I122").Value = 3556.875
$ws.Range("J122").Value = 5633.533
$ws.Range("K122").Value = 10670.625
$ws.Range("L122").Value = 16900.599
$ws.Range("M122").Value = -8220.625
$ws.Range("N122").Value = -21800.599

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 51499.332
$ws.Range("J36").Value = 51499.332
$ws.Range("L36").Value = 51499.332
$ws.Range("N36").Value = -52623.332
$ws.Range("H61").Value = 1184.381
$ws.Range("I61").Value = 1093.6
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 1093.6
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -891.5999999999999
$ws.Range("N61").Value = -3404
$ws.Range("H93").Value = 10872726
$ws.Range("I93").Value = 3493.6428
$ws.Range("J93").Value = 27780420
$ws.Range("K93").Value = 3493.6428
$ws.Range("L93").Value = 27780420
$ws.Range("M93").Value = -2245.6428
$ws.Range("N93").Value = -27782916
$ws.Range("H99").Value = 17722
$ws.Range("I99").Value = 17722
$ws.Range("K99").Value = 17722
$ws.Range("M99").Value = -14727
$ws.Range("H113").Value = 1184.381
$ws.Range("I113").Value = 1093.6
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 1093.6
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 1076.4
$ws.Range("N113").Value = -7340

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 58953.223
$ws.Range("I96").Value = 73125.57000000001
$ws.Range("J96").Value = 9350
$ws.Range("K96").Value = 9350
$ws.Range("L96").Value = 9350
$ws.Range("M96").Value = -71752.57000000001
$ws.Range("N96").Value = -12096

# ---- Special case: CUL row 126, N126 cell removed entirely ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("N126").ClearContents()
